# "some data for sunrise/set for photoresistor"
# Updates timestamp / reading pairs on the "Data In" sheet with newly
# captured sensor data, and refreshes the saved selection / window state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data In")
$ws.Activate()

# xlPasteFormats
$xlPasteFormats = -4122

# Scratch cell used to round-trip a cell's formatting across a .Value
# assignment (plain numeric assignment otherwise clears the cell's
# "quote prefix" flavoured style variant).
$scratch = $ws.Cells.Item(200, 200)

# New Time (col A) / Reading (col B) pairs for each data row
$rows = @(
    @{ Row = 5;  A = 43939.891462662039; B = 3 },
    @{ Row = 8;  A = 43939.890489988429; B = 4 },
    @{ Row = 9;  A = 43939.890559444444; B = 4 },
    @{ Row = 10; A = 43939.890628946756; B = 3 },
    @{ Row = 11; A = 43939.890698425923; B = 4 },
    @{ Row = 12; A = 43939.890767881945; B = 4 },
    @{ Row = 13; A = 43939.890837372688; B = 3 },
    @{ Row = 14; A = 43939.890906828703; B = 3 },
    @{ Row = 15; A = 43939.890976331022; B = 4 },
    @{ Row = 16; A = 43939.891045810182; B = 4 },
    @{ Row = 17; A = 43939.89111527778;  B = 3 },
    @{ Row = 18; A = 43939.891184768516; B = 2 },
    @{ Row = 19; A = 43939.891254259259; B = 3 },
    @{ Row = 20; A = 43939.891323703705; B = 3 },
    @{ Row = 21; A = 43939.891393206017; B = 3 },
    @{ Row = 22; A = 43939.891462662039; B = 3 }
)

foreach ($r in $rows) {
    $timeCell = $ws.Cells.Item($r.Row, 1)
    $readCell = $ws.Cells.Item($r.Row, 2)

    $timeCell.Value = $r.A

    # Preserve the reading cell's existing style (it carries a
    # quote-prefix flagged number format) across the value write.
    $readCell.Copy() | Out-Null
    $scratch.PasteSpecial($xlPasteFormats) | Out-Null

    $readCell.Value = $r.B

    $scratch.Copy() | Out-Null
    $readCell.PasteSpecial($xlPasteFormats) | Out-Null
}

$scratch.Clear() | Out-Null
$excel.CutCopyMode = $false

# Refresh the saved selection on the "Data In" sheet
$ws.Range("A22").Select()

# Best-effort: record the workbook window's last on-screen position
$excel.ActiveWindow.Left = 1170
$excel.ActiveWindow.Top = 1170
